# PluginLive Automation.xlsx - "changes made second commit"
#
# 1. Rename "Sheet1" -> "CTCandJobLocation" (content unchanged).
# 2. Insert two brand-new worksheets, "Questionnaire" and "InterviewWorkflow",
#    between "CTCandJobLocation" and "EligibilityCriteria".
# 3. Expand the existing "EligibilityCriteria" sheet with five new
#    "yearN" columns (I:M) and swap the degree value in A2.

$wb = $excel.ActiveWorkbook

# --- 1. Rename Sheet1 -> CTCandJobLocation -------------------------------
$wsCtc = $wb.Worksheets.Item("Sheet1")
$wsCtc.Name = "CTCandJobLocation"

# --- 2. Insert the two new worksheets (order matters for sheetId) -------
# Questionnaire is created first so it keeps the lower sheetId (5),
# InterviewWorkflow is created second (sheetId 6); both are inserted
# immediately after CTCandJobLocation so the final left-to-right order is
# CorporateLogin, JobDetails, CTCandJobLocation, Questionnaire,
# InterviewWorkflow, EligibilityCriteria.
$wsQuestionnaire = $wb.Worksheets.Add($null, $wsCtc)
$wsQuestionnaire.Name = "Questionnaire"

$wsInterview = $wb.Worksheets.Add($null, $wsQuestionnaire)
$wsInterview.Name = "InterviewWorkflow"

# Re-fetch the EligibilityCriteria handle *after* the inserts above: the
# worksheet collection shifted, and a handle obtained before the inserts
# keeps a stale position.
$wsElig = $wb.Worksheets.Item("EligibilityCriteria")

# --- Questionnaire sheet content -----------------------------------------
$wsQuestionnaire.Range("A1").Value = "question"
$wsQuestionnaire.Range("B1").Value = "option1"
$wsQuestionnaire.Range("C1").Value = "option2"
$wsQuestionnaire.Range("A2").Value = "Are you Okay to Relocate to Chennai?"
$wsQuestionnaire.Range("B2").Value = "Yes"
$wsQuestionnaire.Range("C2").Value = "No"

# Reuse the existing "Consolas" highlight style (style index 3 in
# styles.xml, as already used on JobDetails!C2) for the question cell.
$wb.Worksheets.Item("JobDetails").Range("C2").Copy()
$wsQuestionnaire.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- InterviewWorkflow sheet content --------------------------------------
$wsInterview.Range("A1").Value = "Roundname1"
$wsInterview.Range("B1").Value = "Round1Description"
$wsInterview.Range("C1").Value = "Roundname2"
$wsInterview.Range("D1").Value = "Round2Description"
$wsInterview.Range("E1").Value = "Roundname3"
$wsInterview.Range("F1").Value = "Round3Description"
$wsInterview.Range("G1").Value = "clgName"

$wsInterview.Range("A2").Value = "Assessment"
$wsInterview.Range("B2").Value = "All The Best Do Well"
$wsInterview.Range("C2").Value = "Technical Round"
$wsInterview.Range("D2").Value = "Must Have knowledge regarding the job position"
$wsInterview.Range("E2").Value = "HR Discussion"
$wsInterview.Range("F2").Value = "Discussion regarding the Location, Shift and Package also About the Company"
$wsInterview.Range("G2").Value = "Demo College of engineering"

# --- 3. Expand EligibilityCriteria ---------------------------------------
# Swap the old "Bachelor Of Computer Applications" degree value for the
# new one.
$wsElig.Range("A2").Value = "Bachelor Of Engineeringg"

# Add the five new "yearN" columns with literal year numbers underneath.
$wsElig.Range("I1").Value = "year1"
$wsElig.Range("J1").Value = "year2"
$wsElig.Range("K1").Value = "year3"
$wsElig.Range("L1").Value = "year4"
$wsElig.Range("M1").Value = "year5"

$wsElig.Range("I2").Value = 2020
$wsElig.Range("J2").Value = 2021
$wsElig.Range("K2").Value = 2022
$wsElig.Range("L2").Value = 2023
$wsElig.Range("M2").Value = 2024

# --- Final selections / active sheet, matching the authored workbook -----
$wsQuestionnaire.Range("E3").Select()
$wsInterview.Range("B2").Select()
$wsElig.Range("B2").Select()
$wsCtc.Range("G19").Select()

$wsInterview.Activate()
